$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-27 Monday" "2024-05-28 Tuesday"

Replace-Text "25÷3=" "78÷9="
Replace-Text "63÷6=" "98÷3="
Replace-Text "64÷5=" "24÷3="
Replace-Text "16÷9=" "18÷3="
Replace-Text "87÷4=" "36÷3="

Replace-Text "55÷5=" "57÷6="
Replace-Text "92÷2=" "55÷6="
Replace-Text "97÷7=" "31÷4="
Replace-Text "64÷9=" "23÷9="
Replace-Text "41÷3=" "86÷8="

Replace-Text "94÷6=" "12÷8="
Replace-Text "91÷2=" "67÷7="
Replace-Text "34÷4=" "46÷3="
Replace-Text "89÷9=" "27÷6="
Replace-Text "27÷9=" "33÷7="

Replace-Text "65÷8=" "38÷2="
Replace-Text "71÷9=" "67÷9="
Replace-Text "12÷6=" "66÷4="
Replace-Text "70÷2=" "29÷5="
Replace-Text "93÷5=" "42÷4="

Replace-Text "74÷5=" "66÷2="
Replace-Text "48÷3=" "99÷4="
Replace-Text "38÷4=" "46÷4="
Replace-Text "85÷2=" "18÷3="
Replace-Text "35÷9=" "15÷3="
